# Updated cryptos list on Mon Nov 18 05:12:36 UTC 2024 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for each coin row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.567.40"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "3.104.02"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'242.66"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'625.16"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("D7").Value = "'1.15"
$ws.Range("E7").Value = "  +9.83%  "
$ws.Range("D8").Value = "'0.373"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.745"
$ws.Range("E10").Value = "  +3.92%  "
$ws.Range("E11").Value = "  -18.64%  "
$ws.Range("E12").Value = "  +3.67%  "
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").Value = "'35.28"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "'5.49"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "90.401.67"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "3.674.72"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "3.106.50"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'3.83"
$ws.Range("E19").Value = "  +3.42%  "
$ws.Range("D20").Value = "'14.29"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "'5.78"
$ws.Range("E22").Value = "  +7.42%  "
$ws.Range("D23").Value = "'445.28"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "'5.88"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'92.86"
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "3.262.28"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "'0.178"
$ws.Range("E30").Value = "  +11.15%  "
$ws.Range("D31").Value = "'0.221"
$ws.Range("E31").Value = "  +12.64%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "'4.34"
$ws.Range("E34").Value = "  +34.33%  "
$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  +25.05%  "
$ws.Range("D36").Value = "'26.56"
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").Value = "'7.62"
$ws.Range("E38").Value = "  +8.26%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").Value = "'492.21"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("D41").Value = "'3.59"
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'159.07"
$ws.Range("E46").Value = "  +6.71%  "
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").Value = "'4.56"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "'45.03"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "'1.34"
$ws.Range("E51").Value = "  -1.18%  "

# The leading apostrophe above forces Excel to keep these numeric-looking
# strings as text instead of coercing them to floating-point numbers; reset
# each cell back to the default Normal style afterwards so no stray
# quote-prefix formatting is left behind.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
